# Chapter2-styles.docx: tweak paragraph spacing on a few styles.
#   - Normal:        add spacing-after of 360 twips (18 pt)
#   - Image Caption:  add spacing-after of 480 twips (24 pt) and turn on
#                      "don't add space between paragraphs of the same style"
#                      (this is the w:contextualSpacing toggle)
#   - Figure:         add spacing-after of 480 twips (24 pt)
$d = $word.ActiveDocument

$normal = $d.Styles("Normal")
$normal.ParagraphFormat.SpaceAfter = 18

$imageCaption = $d.Styles("ImageCaption")
$imageCaption.ParagraphFormat.SpaceAfter = 24
$imageCaption.NoSpaceBetweenParagraphsOfSameStyle = $true

$figure = $d.Styles("Figure")
$figure.ParagraphFormat.SpaceAfter = 24
